$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ------------------------------------------------------------------
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Split the mailing-address paragraph (the one in the body, not the
#    one inside the details table) into two paragraphs:
#       "2958 DEVELOPER"
#       "SANTA CLARA, CA 95051"
# ------------------------------------------------------------------
$addrPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "2958 DEVELOPER, SANTA CLARA CA 95051`r" -and
        $cand.Range.Information(12) -eq $false) {
        $addrPara = $cand
        break
    }
}

$addrRange = $addrPara.Range
$addrRange.Find.Execute("2958 DEVELOPER, SANTA CLARA CA 95051", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "2958 DEVELOPER^pSANTA CLARA, CA 95051", 2) | Out-Null

# Fix up the run formatting of the newly created second paragraph so it
# matches the rest of the letter (Arial 11pt) instead of inheriting the
# document default font.
$newAddrPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "SANTA CLARA, CA 95051`r" -and
        $cand.Range.Information(12) -eq $false) {
        $newAddrPara = $cand
        break
    }
}
$newAddrRange = $newAddrPara.Range
$newAddrRange.Font.Name = "Arial"
$newAddrRange.Font.NameBi = "Arial"
$newAddrRange.Font.Size = 11
$newAddrRange.Font.SizeBi = 11

# ------------------------------------------------------------------
# 3. Remove the empty "No Spacing" paragraph right after the
#    "... Board of Directors" signature line.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -match "Board of Directors") {
        $nextPara = $d.Paragraphs.Item($i + 1)
        if ($nextPara.Range.Text -eq "`r" -and $nextPara.Style.NameLocal -eq "No Spacing") {
            $nextPara.Range.Delete()
        }
        break
    }
}
